$wb = $excel.ActiveWorkbook

# --- Hoja1: update the free-text "Conversión del día" note in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = $wsHoja1.Range("A1").Value()
$text = $text -replace [regex]::Escape("1000 Bs = 12.39 = 50000.0 pesos"), "1000 Bs = 12.48 = 50489.28 pesos"
$text = $text -replace [regex]::Escape("50000.0 pesos = 12.38 = 975.57 Bs"), "50489.28 pesos = 12.41 = 967.61 Bs"
$wsHoja1.Range("A1").Value = $text

# --- tasas: update the transfi rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 80.15900000000001
$wsTasas.Range("O10").Value = 4047.17
$wsTasas.Range("N12").Value = 4069.98
$wsTasas.Range("O12").Value = 78
